$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Map of row number -> @(DAMSLTag, DialogAct)
$changes = @{
    3   = @("sd", "Statement-non-opinion")
    25  = @("sv", "Statement-opinion")
    26  = @("sv", "Statement-opinion")
    31  = @("sv", "Statement-opinion")
    32  = @("aa", "Agree/Accept")
    43  = @("aa", "Agree/Accept")
    53  = @("aa", "Agree/Accept")
    62  = @("aa", "Agree/Accept")
    63  = @("aa", "Agree/Accept")
    65  = @("ba", "Appreciation")
    68  = @("sd", "Statement-non-opinion")
    69  = @("aa", "Agree/Accept")
    71  = @("sv", "Statement-opinion")
    74  = @("%", "Uninterpretable")
    83  = @("sd", "Statement-non-opinion")
    94  = @("b", "Acknowledge (Backchannel)")
    101 = @("b", "Acknowledge (Backchannel)")
    107 = @("aa", "Agree/Accept")
    111 = @("ba", "Appreciation")
    117 = @("sv", "Statement-opinion")
    132 = @("sd", "Statement-non-opinion")
    136 = @("%", "Uninterpretable")
}

foreach ($row in $changes.Keys) {
    $values = $changes[$row]
    $ws.Cells.Item($row, 9).Value = $values[0]
    $ws.Cells.Item($row, 10).Value = $values[1]
}
